# The deck carries two DrawingML theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Office"     (theme name "Office Theme")
#   ppt/theme/theme2.xml -> clrScheme "Red Violet"  (theme name "Integral")
# theme2.xml is the live theme (wired to the one slide master / the
# presentation itself), so "Integral" is what actually renders today.
#
# The authored change swaps the two parts' contents, which makes the
# presentation's live theme become the "Office Theme" palette instead of
# "Integral". Reproduce that effect on the live theme by rewriting its 12
# theme colors - in clrScheme order dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink - to the "Office Theme" values.

function HexToRgb([string]$hex) {
    # PowerPoint's RGB property packs colors as 0x00BBGGRR (red in the low
    # byte), matching the native RGB() layout used throughout the object
    # model - NOT the 0xRRGGBB order the OOXML <a:srgbClr val="RRGGBB"/>
    # attribute uses.
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $cs.Item($i).RGB = HexToRgb $officeThemeColors[$i - 1]
}
